$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Row 2 ("Yes" row), column 2 (QRPs list) ---
$cell1 = $t.Cell(2, 2).Range

$cell1.Find.Execute(
    "Choosing a poor model specification, ", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Choosing a poor model specification, Choosing biased measurements, Choosing overlapping measures to find significant results, ",
    2) | Out-Null

$cell1.Find.Execute(
    "Using biased measurements, Using irrelevant references, Using measurement overlap to find significant results, Using unjustified references",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Using irrelevant references, Using unjustified references",
    2) | Out-Null

# --- Row 3 ("Maybe" row), column 2 (QRPs list) ---
$cell2 = $t.Cell(3, 2).Range

$cell2.Find.Execute(
    "Creating multiple publications", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Choosing biased manipulations, Creating multiple publications",
    2) | Out-Null

$cell2.Find.Execute(
    "Modifying measurements, Not linking the Preregistration", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Mixing pilot and main study data, Modifying measurements, Not linking the preregistration",
    2) | Out-Null

$cell2.Find.Execute(
    "Redefining group membership rules, Retaining pilot data, Selecting", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Redefining group membership rules, Selecting",
    2) | Out-Null

$cell2.Find.Execute(
    "Using ad hoc exclusion criteria for participants, Using biased manipulations, Variable transformation fishing",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Using ad hoc exclusion criteria for participants, Variable transformation fishing",
    2) | Out-Null

Write-Output "Cell(2,2): $($t.Cell(2,2).Range.Text)"
Write-Output "Cell(3,2): $($t.Cell(3,2).Range.Text)"
